$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 498, shifting the existing rows (498..523) down to (499..524)
$ws.Rows(498).Insert()

# Populate the newly inserted row 498 with the new weekly price record
$ws.Cells.Item(498, 1).Value  = 3
$ws.Cells.Item(498, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(498, 3).Value  = "Coquimbo"
$ws.Cells.Item(498, 4).Value  = 45041
$ws.Cells.Item(498, 5).Value  = 5
$ws.Cells.Item(498, 6).Value  = 100112012
$ws.Cells.Item(498, 7).Value  = "Espinaca"
$ws.Cells.Item(498, 8).Value  = "Sin especificar"
$ws.Cells.Item(498, 9).Value  = "Primera"
$ws.Cells.Item(498, 10).Value = 80
$ws.Cells.Item(498, 11).Value = 4500
$ws.Cells.Item(498, 12).Value = 4500
$ws.Cells.Item(498, 13).Value = 4500
$ws.Cells.Item(498, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(498, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(498, 16).Value = 1500
$ws.Cells.Item(498, 17).Value = 3
$ws.Cells.Item(498, 18).Value = "Hortaliza"
